# The commit inserts one new blank row above the existing data on Sheet1,
# pushing every row down by one (old row 1 -> new row 2, ..., old row 68 -> new
# row 69). This matches the diff: dimension goes from A1:C68 to A2:C69, and
# every <row>/<c> value shifts down by exactly one row (cell styles, like the
# s="1" font on the old B36/new B37, travel with their row).
#
# Reproduce it the way a user would in the UI: click the row-1 header to
# select the whole row, then insert a new row above it (shifting everything
# down).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Rows(1).Select()
$ws.Rows(1).Insert()
